$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Columns("E").Delete()
$ws2.Columns("D").Delete()

$ws1.Activate() | Out-Null
$ws1.Range("C15").Select() | Out-Null
